$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 460.14285
$ws.Range("J8").Value = 19
$ws.Range("L8").Value = 57
$ws.Range("N8").Value = -335
$ws.Range("H9").Value = 219.66667
$ws.Range("J9").Value = 347.5
$ws.Range("L9").Value = 347.5
$ws.Range("N9").Value = -685.5
$ws.Range("H41").Value = 2600.3
$ws.Range("I41").Value = 2478.7144
$ws.Range("K41").Value = 2478.7144
$ws.Range("M41").Value = -2038.7144
$ws.Range("H103").Value = 958
$ws.Range("J103").Value = 1000
$ws.Range("L103").Value = 3000
$ws.Range("N103").Value = -4172
$ws.Range("H112").Value = 5716.5713
$ws.Range("J112").Value = 5956.4907
$ws.Range("L112").Value = 17869.4721
$ws.Range("N112").Value = -20085.4721
$ws.Range("H127").Value = 2223
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H137").Value = 19610372
$ws.Range("I137").Value = 40001996
$ws.Range("J137").Value = 3041.9614
$ws.Range("K137").Value = 120005988
$ws.Range("L137").Value = 9125.8842
$ws.Range("M137").Value = -120003438
$ws.Range("N137").Value = -14225.8842

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 2037.5
$ws.Range("I25").Value = 575
$ws.Range("K25").Value = 575
$ws.Range("M25").Value = -173
$ws.Range("H45").Value = 2647.3914
$ws.Range("I45").Value = 1888.9474
$ws.Range("K45").Value = 1888.9474
$ws.Range("M45").Value = -1511.9474
$ws.Range("H74").Value = 3366.8
$ws.Range("I74").Value = 944.6667
$ws.Range("K74").Value = 944.6667
$ws.Range("M74").Value = -70.66669999999999
$ws.Range("H77").Value = 3366.8
$ws.Range("I77").Value = 944.6667
$ws.Range("K77").Value = 4723.3335
$ws.Range("M77").Value = -355.3334999999997
$ws.Range("H82").Value = 50181
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 50181
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 50181
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -50903
$ws.Range("H85").Value = 50181
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 50181
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 50181
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -52677
$ws.Range("H102").Value = 37344.8
$ws.Range("I102").Value = 52293.43
$ws.Range("K102").Value = 52293.43
$ws.Range("M102").Value = -50671.43

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1964.174
$ws.Range("I105").Value = 1979.1
$ws.Range("J105").Value = 1864.6666
$ws.Range("K105").Value = 1979.1
$ws.Range("L105").Value = 1864.6666
$ws.Range("M105").Value = -232.0999999999999
$ws.Range("N105").Value = -5358.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20868.959
$ws.Range("I31").Value = 7113.9375
$ws.Range("K31").Value = 7113.9375
$ws.Range("M31").Value = -6818.9375
$ws.Range("H34").Value = 20868.959
$ws.Range("I34").Value = 7113.9375
$ws.Range("K34").Value = 7113.9375
$ws.Range("M34").Value = -6911.9375
$ws.Range("H58").Value = 689948.5600000001
$ws.Range("I58").Value = 825937.4
$ws.Range("K58").Value = 825937.4
$ws.Range("M58").Value = -825734.4
$ws.Range("H99").Value = 3220.75
$ws.Range("I99").Value = 2579.4
$ws.Range("K99").Value = 2579.4
$ws.Range("M99").Value = -1081.4
$ws.Range("H126").Value = 3220.75
$ws.Range("I126").Value = 2579.4
$ws.Range("K126").Value = 7738.200000000001
$ws.Range("M126").Value = -5268.200000000001
$ws.Range("H132").Value = 38657030
$ws.Range("I132").Value = 55557500
$ws.Range("K132").Value = 166672500
$ws.Range("M132").Value = -166669970
$ws.Range("H136").Value = 689948.5600000001
$ws.Range("I136").Value = 825937.4
$ws.Range("K136").Value = 2477812.2
$ws.Range("M136").Value = -2475262.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 10665.5
$ws.Range("I3").Value = 11110
$ws.Range("J3").Value = 9998.75
$ws.Range("K3").Value = 33330
$ws.Range("L3").Value = 29996.25
$ws.Range("M3").Value = -33218
$ws.Range("N3").Value = -30220.25
$ws.Range("H4").Value = 33108446
$ws.Range("I4").Value = 40261500
$ws.Range("K4").Value = 120784500
$ws.Range("M4").Value = -120784388
$ws.Range("H22").Value = 862.25
$ws.Range("J22").Value = 2999
$ws.Range("L22").Value = 8997
$ws.Range("N22").Value = -9335
$ws.Range("H27").Value = 862.25
$ws.Range("J27").Value = 2999
$ws.Range("L27").Value = 8997
$ws.Range("N27").Value = -9201
$ws.Range("H41").Value = 112.5
$ws.Range("J41").Value = 150
$ws.Range("L41").Value = 450
$ws.Range("N41").Value = -1126
$ws.Range("H42").Value = 7500
$ws.Range("J42").Value = 7500
$ws.Range("L42").Value = 22500
$ws.Range("N42").Value = -23568
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H48").Value = 4062.5
$ws.Range("I48").Value = 1500
$ws.Range("K48").Value = 4500
$ws.Range("M48").Value = -4250
$ws.Range("H51").Value = 4625
$ws.Range("I51").Value = 3500
$ws.Range("K51").Value = 10500
$ws.Range("M51").Value = -10040
$ws.Range("H54").Value = 5475
$ws.Range("J54").Value = 5475
$ws.Range("L54").Value = 16425
$ws.Range("N54").Value = -17543
$ws.Range("H55").Value = 1566
$ws.Range("I55").Value = 1254.6666
$ws.Range("J55").Value = 2500
$ws.Range("K55").Value = 3763.9998
$ws.Range("L55").Value = 7500
$ws.Range("M55").Value = -3586.9998
$ws.Range("N55").Value = -7854
$ws.Range("H57").Value = 3320.5454
$ws.Range("I57").Value = 1417
$ws.Range("K57").Value = 4251
$ws.Range("M57").Value = -3692
$ws.Range("H74").Value = 15561.667
$ws.Range("J74").Value = 15561.667
$ws.Range("L74").Value = 46685.001
$ws.Range("N74").Value = -48807.001
$ws.Range("H77").Value = 15561.667
$ws.Range("J77").Value = 15561.667
$ws.Range("L77").Value = 140055.003
$ws.Range("N77").Value = -150663.003
$ws.Range("H131").Value = 9367.556
$ws.Range("I131").Value = 835.05884
$ws.Range("J131").Value = 23872.8
$ws.Range("K131").Value = 2505.17652
$ws.Range("L131").Value = 71618.39999999999
$ws.Range("M131").Value = 2534.82348
$ws.Range("N131").Value = -81698.39999999999
$ws.Range("H132").Value = 950
$ws.Range("I132").Value = 900
$ws.Range("K132").Value = 8100
$ws.Range("M132").Value = -5570

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2832.0908
$ws.Range("I102").Value = 2215.125
$ws.Range("K102").Value = 2215.125
$ws.Range("M102").Value = -593.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3924.2
$ws.Range("I40").Value = 3873.6667
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 3873.6667
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -3737.6667
$ws.Range("N40").Value = -4272
$ws.Range("H68").Value = 2428.1428
$ws.Range("I68").Value = 2399.4
$ws.Range("K68").Value = 2399.4
$ws.Range("M68").Value = -1650.4
$ws.Range("H71").Value = 2428.1428
$ws.Range("I71").Value = 2399.4
$ws.Range("K71").Value = 11997
$ws.Range("M71").Value = -8253
$ws.Range("H136").Value = 3857.5
$ws.Range("I136").Value = 2884.4546
$ws.Range("K136").Value = 8653.363799999999
$ws.Range("M136").Value = -6103.363799999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2740
$ws.Range("I62").Value = 2166.6667
$ws.Range("K62").Value = 2166.6667
$ws.Range("M62").Value = -1542.6667
$ws.Range("H65").Value = 2740
$ws.Range("I65").Value = 2166.6667
$ws.Range("K65").Value = 10833.3335
$ws.Range("M65").Value = -7713.333500000001
$ws.Range("H100").Value = 1570.375
$ws.Range("I100").Value = 843.5
$ws.Range("J100").Value = 2781.8333
$ws.Range("K100").Value = 1687
$ws.Range("L100").Value = 5563.6666
$ws.Range("M100").Value = -1146
$ws.Range("N100").Value = -6645.6666
$ws.Range("H122").Value = 2334.3914
$ws.Range("I122").Value = 2142.7856
$ws.Range("K122").Value = 6428.3568
$ws.Range("M122").Value = -3978.3568
$ws.Range("H132").Value = 26194116
$ws.Range("I132").Value = 28208124
$ws.Range("K132").Value = 84624372
$ws.Range("M132").Value = -84621842
